$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text parses as a plain number need an explicit
# Text number format first, otherwise Excel auto-converts the assignment
# (e.g. "189.40" -> 189.4) and the trailing zero / exact text is lost.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D14", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D28", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '68.858.46'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '3.334.11'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '189.40'
$ws.Range("E5").Value = '  +1.90%  '
$ws.Range("D6").Value = '588.10'
$ws.Range("E6").Value = '  +1.28%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("D9").Value = '0.132'
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("D10").Value = '6.74'
$ws.Range("E10").Value = '  +2.49%  '
$ws.Range("D11").Value = '0.414'
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("D12").Value = '3.923.86'
$ws.Range("E12").Value = '  +2.50%  '
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("D14").Value = '28.11'
$ws.Range("E14").Value = '  +2.27%  '
$ws.Range("D15").Value = '68.995.90'
$ws.Range("E15").Value = '  +1.51%  '
$ws.Range("D16").Value = '0.0000169'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").Value = '3.327.52'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = '445.59'
$ws.Range("E18").Value = '  +11.71%  '
$ws.Range("D19").Value = '5.79'
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("D20").Value = '13.70'
$ws.Range("E20").Value = '  +1.23%  '
$ws.Range("D21").Value = '7.80'
$ws.Range("E21").Value = '  +2.47%  '
$ws.Range("D22").Value = '75.53'
$ws.Range("E22").Value = '  +5.96%  '
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").Value = '3.506.24'
$ws.Range("E24").Value = '  +2.65%  '
$ws.Range("D25").Value = '0.522'
$ws.Range("E25").Value = '  +2.01%  '
$ws.Range("D26").Value = '0.0000120'
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("E27").Value = '  +1.25%  '
$ws.Range("D28").Value = '9.37'
$ws.Range("E28").Value = '  -1.40%  '
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D30").Value = '2.00'
$ws.Range("E30").Value = '  +2.51%  '
$ws.Range("D31").Value = '23.20'
$ws.Range("E31").Value = '  +2.35%  '
$ws.Range("D32").Value = '5.49'
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("D33").Value = '1.27'
$ws.Range("E33").Value = '  +1.30%  '
$ws.Range("D34").Value = '6.93'
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = '1.55'
$ws.Range("E36").Value = '  +5.78%  '
$ws.Range("D37").Value = '163.64'
$ws.Range("E37").Value = '  +0.63%  '
$ws.Range("D38").Value = '1.92'
$ws.Range("E38").Value = '  +1.50%  '
$ws.Range("D39").Value = '27.13'
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("D40").Value = '4.57'
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("D41").Value = '0.799'
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("D42").Value = '6.45'
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("D43").Value = '2.700.17'
$ws.Range("E43").Value = '  +1.00%  '
$ws.Range("D44").Value = '2.48'
$ws.Range("E44").Value = '  +1.76%  '
$ws.Range("D45").Value = '41.22'
$ws.Range("E45").Value = '  +1.21%  '
$ws.Range("D46").Value = '0.0684'
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("D47").Value = '25.18'
$ws.Range("E47").Value = '  +1.94%  '
$ws.Range("D48").Value = '328.38'
$ws.Range("E48").Value = '  -1.85%  '
$ws.Range("D49").Value = '0.0282'
$ws.Range("E49").Value = '  +2.51%  '
$ws.Range("D50").Value = '32.30'
$ws.Range("E50").Value = '  +5.29%  '
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +3.08%  '
